$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking Price values in column D stay as text (matching the
# original inline-string representation) instead of being auto-converted to numbers.
$textCells = @("D4","D5","D6","D8","D9","D10","D11","D12","D14","D15","D17","D18","D19","D20","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the refreshed coin data (price + 1h volume, and for rows that shifted up
# one position in the ranking, also the coin name + link).
$ws.Range("D2").Value = '27.277.18'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '1.863.61'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("D4").Value = '1.021'
$ws.Range("E4").Value = '  +1.43%  '
$ws.Range("D5").Value = '313.03'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = '1.019'
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  +1.79%  '
$ws.Range("D8").Value = '0.3737'
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").Value = '0.07445'
$ws.Range("E9").Value = '  +4.21%  '
$ws.Range("D10").Value = '0.9373'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").Value = '20.68'
$ws.Range("E11").Value = '  +5.95%  '
$ws.Range("D12").Value = '0.07882'
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").Value = '1.872.90'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '5.434'
$ws.Range("E14").Value = '  +2.93%  '
$ws.Range("D15").Value = '6.543'
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").Value = '1.022'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '0.000008792'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").Value = '1.019'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '14.82'
$ws.Range("E20").Value = '  +2.48%  '
$ws.Range("D21").Value = '27.306.45'
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").Value = '5.121'
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("D23").Value = '10.69'
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '1.956'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '154.22'
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("E26").Value = '  +2.00%  '
$ws.Range("D27").Value = '2.007'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").Value = '116.06'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("D29").Value = '4.997'
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("D30").Value = '0.08917'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").Value = '3.348'
$ws.Range("E31").Value = '  +4.25%  '
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").Value = '4.570'
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("D34").Value = '0.7456'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = '2.677'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '0.02054'
$ws.Range("E36").Value = '  +5.81%  '
$ws.Range("D37").Value = '1.126'
$ws.Range("E37").Value = '  +3.52%  '
$ws.Range("D38").Value = '0.05288'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.5363'
$ws.Range("E39").Value = '  +3.40%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '7.128'
$ws.Range("E40").Value = '  +2.39%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.1537'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '8.391'
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '10.63'
$ws.Range("E43").Value = '  +2.25%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.4832'
$ws.Range("E44").Value = '  +2.94%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.020'
$ws.Range("E45").Value = '  +1.38%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '1.669'
$ws.Range("E46").Value = '  +4.87%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '103.04'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '66.72'
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06084'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '0.9012'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '36.84'
$ws.Range("E51").Value = '  +1.75%  '
